# Review document cleanup (Sprint 3 / ReviewDoc.xlsx)
#   - Highlight the header row (A7) with the same "done section" fill used
#     elsewhere in the sheet (no border).
#   - Fill in the previously-empty rows 29-33 with the Transactie/Asset
#     classes that were reviewed, each with their own "maker" note.
#   - Move the active selection down to A34 (next empty row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- A7: restyle to the "section divider" fill (fill only, no border) ---
$fillOnlySource = $ws.Range("A9")
$fillOnlySource.Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Borders.LineStyle = -4142

# Source cell carrying the highlighted "class name" style (fill + border)
# used by the other populated rows (e.g. A27, A28) in this table.
$classNameStyleSource = $ws.Range("A27")

# --- Row 29: TransactiePaginaDto / Carmen ---
$classNameStyleSource.Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = "TransactiePaginaDto"
$ws.Range("B29").Value = "Carmen"

# --- Row 30: TransactieStartDto / Carmen ---
$classNameStyleSource.Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A30").Value = "TransactieStartDto"
$ws.Range("B30").Value = "Carmen"

# --- Row 31: TransactieService / Carmen en (style unchanged) ---
$ws.Range("A31").Value = "TransactieService"
$ws.Range("B31").Value = "Carmen en"

# --- Row 32: AssetDto / Carmen ---
$classNameStyleSource.Copy()
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A32").Value = "AssetDto"
$ws.Range("B32").Value = "Carmen"

# --- Row 33: AssetNotExistsException / Carmen ---
$classNameStyleSource.Copy()
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A33").Value = "AssetNotExistsException"
$ws.Range("B33").Value = "Carmen"

# --- Move selection to the next empty row ---
$null = $ws.Range("A34").Select()

Write-Host "ReviewDoc cleanup applied"
